$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 431.66666
$ws.Range("I2").Value = 597.5
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 597.5
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -484.5
$ws.Range("N2").Value = -326
$ws.Range("H32").Value = 3178.5
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3178.5
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3178.5
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -3830.5
$ws.Range("H51").Value = 26889.2
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 26889.2
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 26889.2
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -27857.2
$ws.Range("H129").Value = 871.39343
$ws.Range("I129").Value = 483.33334
$ws.Range("K129").Value = 1450.00002
$ws.Range("M129").Value = 3549.99998
$ws.Range("H137").Value = 2441301.5
$ws.Range("I137").Value = 4001837.5
$ws.Range("J137").Value = 2964.1875
$ws.Range("K137").Value = 12005512.5
$ws.Range("L137").Value = 8892.5625
$ws.Range("M137").Value = -12002962.5
$ws.Range("N137").Value = -13992.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6432538.5
$ws.Range("I32").Value = 7962853.5
$ws.Range("J32").Value = 5213.6
$ws.Range("K32").Value = 7962853.5
$ws.Range("L32").Value = 5213.6
$ws.Range("M32").Value = -7962566.5
$ws.Range("N32").Value = -5787.6
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 400
$ws.Range("K12").Value = 400
$ws.Range("M12").Value = -232
$ws.Range("H29").Value = 1633.3334
$ws.Range("I29").Value = 1633.3334
$ws.Range("K29").Value = 1633.3334
$ws.Range("M29").Value = -1344.3334
$ws.Range("H36").Value = 660
$ws.Range("I36").Value = 500
$ws.Range("J36").Value = 980
$ws.Range("K36").Value = 500
$ws.Range("L36").Value = 980
$ws.Range("M36").Value = 34
$ws.Range("N36").Value = -2048
$ws.Range("H37").Value = 616.6667
$ws.Range("I37").Value = 340
$ws.Range("K37").Value = 340
$ws.Range("M37").Value = -203
$ws.Range("H48").Value = 70000
$ws.Range("J48").Value = 70000
$ws.Range("L48").Value = 70000
$ws.Range("N48").Value = -70830

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1796.9796
$ws.Range("J31").Value = 2706.7058
$ws.Range("L31").Value = 2706.7058
$ws.Range("N31").Value = -3296.7058
$ws.Range("H34").Value = 1796.9796
$ws.Range("J34").Value = 2706.7058
$ws.Range("L34").Value = 2706.7058
$ws.Range("N34").Value = -3110.7058
$ws.Range("H41").Value = 1650
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""
$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -4264
$ws.Range("N51").Value = ""
$ws.Range("H59").Value = 40127
$ws.Range("J59").Value = 40127
$ws.Range("L59").Value = 40127
$ws.Range("N59").Value = -42417
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4652
$ws.Range("N61").Value = ""
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = ""
$ws.Range("H94").Value = 2957.3044
$ws.Range("I94").Value = 15604
$ws.Range("J94").Value = 1060.3
$ws.Range("K94").Value = 15604
$ws.Range("L94").Value = 1060.3
$ws.Range("M94").Value = -15153
$ws.Range("N94").Value = -1962.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1257.3572
$ws.Range("I122").Value = 534.6667
$ws.Range("J122").Value = 1454.4546
$ws.Range("K122").Value = 4812.0003
$ws.Range("L122").Value = 13090.0914
$ws.Range("M122").Value = -2362.0003
$ws.Range("N122").Value = -17990.0914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1634.6471
$ws.Range("I113").Value = 1462.7142
$ws.Range("J113").Value = 2437
$ws.Range("K113").Value = 1462.7142
$ws.Range("L113").Value = 2437
$ws.Range("M113").Value = 707.2858000000001
$ws.Range("N113").Value = -6777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 776.6818
$ws.Range("I22").Value = 328.66666
$ws.Range("J22").Value = 1086.8462
$ws.Range("K22").Value = 328.66666
$ws.Range("L22").Value = 1086.8462
$ws.Range("M22").Value = -33.66665999999998
$ws.Range("N22").Value = -1676.8462
$ws.Range("H27").Value = 776.6818
$ws.Range("I27").Value = 328.66666
$ws.Range("J27").Value = 1086.8462
$ws.Range("K27").Value = 328.66666
$ws.Range("L27").Value = 1086.8462
$ws.Range("M27").Value = -221.66666
$ws.Range("N27").Value = -1300.8462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 11404.333
$ws.Range("I54").Value = 325
$ws.Range("J54").Value = 13620.2
$ws.Range("K54").Value = 325
$ws.Range("L54").Value = 13620.2
$ws.Range("M54").Value = 195
$ws.Range("N54").Value = -14660.2
$ws.Range("H81").Value = 3294.818
$ws.Range("I81").Value = 2138.111
$ws.Range("J81").Value = 8500
$ws.Range("K81").Value = 4276.222
$ws.Range("L81").Value = 17000
$ws.Range("M81").Value = -3215.222
$ws.Range("N81").Value = -19122
$ws.Range("H84").Value = 3294.818
$ws.Range("I84").Value = 2138.111
$ws.Range("J84").Value = 8500
$ws.Range("K84").Value = 21381.11
$ws.Range("L84").Value = 85000
$ws.Range("M84").Value = -16077.11
$ws.Range("N84").Value = -95608
